$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2) | Out-Null
}

# "English" appears twice (hyperlink label + plain run) -- both become الإنجليزية
Replace-Text "English" "الإنجليزية"

# Language list after the hyperlink (keep the leading space out of the Find
# text -- including it causes the run to inherit the adjacent hyperlink's
# formatting in this runtime, so match starting at the "/")
Replace-Text "/ Portuguese / French / Thai / Vietnamese / Spanish" "/البرتغالية/الفرنسية/التايلندية/الفيتنامية/الإسبانية"

# Table cell "Brief" heading
Replace-Text "Brief" "المضمون"

# Table cell description line
Replace-Text "An email sent to the confirmed attendees of the event. It will be sent via customer.io" "An email sent to the confirmed attendees of the event. سيتم إرسالها عبر customer.io"

# Table cell "Target audience" heading
Replace-Text "Target audience" "الجمهور المستهدف"

# Heading2 title
Replace-Text "Travel checklist: here's what you need" "قائمة التحقق الخاصة بالسفر: إليك ما تحتاجه"

# Greeting line
Replace-Text "Hi " "مرحبًا "
Replace-Text "[PARTNER NAME]" "[اسم الشريك]"

# The comma right after [PARTNER NAME] needs to become an Arabic comma.
# Scope the find to that specific paragraph so we don't touch other commas.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*اسم الشريك*") {
        $r = $p.Range
        $r.Find.Execute(",", $true, $false, $false, $false, $false, $true, 1, $false, "،", 2) | Out-Null
        break
    }
}

# Checklist intro line
Replace-Text "Here’s a checklist of the necessary items for your trip: " "فيما يلي قائمة مرجعية بالعناصر الضرورية لرحلتك: "

# Passport bullet
Replace-Text "Passport " "جواز سفر "

# Yellow fever bullet note (only the second sentence is translated)
Replace-Text "For travellers from yellow fever endemic countries, follow the requirements set by your country. Vaccination should be done no less than 14 days prior to the journey. " "For travellers from yellow fever endemic countries, follow the requirements set by your country. يجب أن يتم التطعيم قبل 14 يومًا على الأقل من الرحلة. "

# Travel itinerary bullet
Replace-Text "A digital or printed copy of the travel itinerary" "نسخة رقمية أو مطبوعة من خط سير السفر"

# Smart casual bullet
Replace-Text "Smart casual attire for the conference" "ملابس غير رسمية أنيقة للمؤتمر"

# Black tie bullet
Replace-Text "Black tie attire for the Gala dinner" "ربطة عنق سوداء لحفل العشاء"

# Contact via live chat line
Replace-Text "If you have any questions, please contact us via " "إذا كانت لديك أي أسئلة، فاتصل بنا:  "
Replace-Text "live chat" "الدردشة الحية"

# Contact country manager line
Replace-Text "If you have any questions, please contact your country manager, " "إذا كانت لديك أي أسئلة، فيُرجى الاتصال بمدير بلدك  "

# Comment text "choose either one"
foreach ($c in $d.Comments) {
    if ($c.Range.Text -eq "choose either one") {
        $c.Range.Find.Execute("choose either one", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "اختر أيًا منهما", 2) | Out-Null
    }
}
